$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gen")
$ws.Activate()

$ws.Range("B5").Value = 0.4

$ws.Range("C2:C3").Select()
$ws.Cells.Item(3, 3).Activate()
